# Add a "Spell Icon" and "Overlay" column to both tables (Table1 on sheet
# "Buffs and Utility" and Table2 on sheet "Debuffs and Attacks"), to support
# the new HUD overlay displaying currently slotted spells.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    if ($ws.ListObjects.Count -lt 1) { continue }
    $lo = $ws.ListObjects.Item(1)

    $lastHeaderCell = $lo.HeaderRowRange.Item(1, $lo.ListColumns.Count)
    $lastDataCol = $lastHeaderCell.Column

    # Add the two new trailing columns to the table.
    $iconCol = $lo.ListColumns.Add()
    $ws.Cells.Item(1, $iconCol.Index).Value = "Spell Icon"

    $overlayCol = $lo.ListColumns.Add()
    $ws.Cells.Item(1, $overlayCol.Index).Value = "Overlay"

    # Match the header formatting of the preceding header cell (bold, centered).
    $newHeaderRange = $ws.Range($ws.Cells.Item(1, $iconCol.Index), $ws.Cells.Item(1, $overlayCol.Index))
    $newHeaderRange.Font.Bold = $true
    $newHeaderRange.HorizontalAlignment = -4108

    # Fill every data row with "Yes", matching the other feature columns, and
    # center the values like the rest of the data columns.
    $lastRow = $lo.Range.Row + $lo.Range.Rows.Count - 1
    $firstDataRow = $lo.Range.Row + 1
    if ($lastRow -ge $firstDataRow) {
        $iconData = $ws.Range($ws.Cells.Item($firstDataRow, $iconCol.Index), $ws.Cells.Item($lastRow, $overlayCol.Index))
        $iconData.Value = "Yes"
        $iconData.HorizontalAlignment = -4108
    }

    # Approximate the column widths Excel's AutoFit would have produced.
    $ws.Columns.Item($iconCol.Index).ColumnWidth = 13.25
    $ws.Columns.Item($overlayCol.Index).ColumnWidth = 11.5
}
